# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2304"
#   "<header>_new" -> "<header>_FV2310"
# and turn the data range into a proper Excel Table, with the header row frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1, columns A:U) -------------------------
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -like "*_old") {
            $cell.Value = ($val -replace "_old$", "_FV2304")
        } elseif ($val -like "*_new") {
            $cell.Value = ($val -replace "_new$", "_FV2310")
        }
    }
}

# --- 2. Convert the used range into an Excel Table named "Table1" ----------
$dataRange = $ws.Range("A1:U94")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (split below row 1) ---------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
